$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Main GUPRI (i.e. PID) under which all terms are defined. Preference is to use PURLs or W3IDs as they provide permanent resolvable identifiers."
$ws.Range("E1").ClearContents()
$ws.Range("C3").Value = "http://purl.org/m4m/"
$ws.Range("E3").Value = "Prefix for our controlled vocabulary since it is rather tedious to write long URLs all the time"
$ws.Range("G3").ClearContents()
$ws.Range("C4").Value = "http://www.w3.org/2004/02/skos/core#"
$ws.Range("E4").Value = "Prefix for SKOS Onotlogy is which our base for defining SKOS based controlled vocabulary"
$ws.Range("G4").ClearContents()
$ws.Range("C5").Value = "http://purl.org/pav/ "
$ws.Range("E5").Value = "Prefix for Provenance, Authoring and Versioning Onotlogy which properties such as version and createdOn we will use to describe our controlled vocabulary"
$ws.Range("G5").ClearContents()
$ws.Range("C6").Value = "http://purl.org/dc/terms/ "
$ws.Range("E6").Value = "Prefix for Dublin Core (Terms) Ontology which properties such as title, description, rights, source, etc. we will use to describe our controlled vocabulary as well to define its terms"
$ws.Range("G6").ClearContents()
$ws.Range("C7").Value = "https://w3id.org/iadopt/ont/"
$ws.Range("E7").Value = "I-ADOPT Ontology"
$ws.Range("G7").ClearContents()
$ws.Range("C8").Value = "https://w3id.org/env/puv#"
$ws.Range("E8").Value = "A simple ontology which implements the Parameter Usage `nVocabulary semantic model, as described at `nhttps://github.com/nvs-vocabs/P01."
$ws.Range("G8").ClearContents()
$ws.Range("C9").Value = "http://www.w3.org/ns/sosa/"
$ws.Range("E9").Value = "Semantic Sensor Network Ontology"
$ws.Range("G9").ClearContents()
$ws.Range("C10").Value = "http://www.w3.org/2001/XMLSchema#"
$ws.Range("E10").ClearContents()
$ws.Range("C11").Value = "http://qudt.org/schema/qudt/"
$ws.Range("E11").Value = "Prefix for ""Quantity, Unit, Dimension and Type"" schema that is used to model physical units."
$ws.Range("G11").ClearContents()
$ws.Range("C12").Value = "http://qudt.org/vocab/unit/"
$ws.Range("E12").Value = "Prefix for QUDT Vocabulary of Units which terms we will use to semantically define units of measurements for terms (i.e., variables) we are defining "
$ws.Range("G12").ClearContents()
$ws.Range("C14").Value = "Title of the vocabulary"
$ws.Range("E14").ClearContents()
$ws.Range("C15").Value = "Description of the controlled vocabulary"
$ws.Range("E15").ClearContents()
$ws.Range("C16").Value = "Home page of community creating vocabulary"
$ws.Range("E16").ClearContents()
$ws.Range("C17").Value = "An ORCID ID of the vocabulary creator"
$ws.Range("E17").ClearContents()
$ws.Range("C18").Value = "License under which the vocabulary is provided"
$ws.Range("E18").ClearContents()
$ws.Range("C19").Value = "Vocabulary version"
$ws.Range("E19").ClearContents()
$ws.Range("C20").Value = "Date when vocabulary was initially created"
$ws.Range("E20").ClearContents()
$ws.Range("B21").Value = "2023-09-13T15:01:00+00:00"
$ws.Range("C21").Value = "Automatic update when vocabulary is updated"
$ws.Range("E21").ClearContents()
$ws.Range("AA23").ClearContents()
$ws.Range("E23").Value = "skos:definition@en"
$ws.Range("F23").Value = "dct:source(separator="","")"
$ws.Range("G23").Value = "skos:broader(separator="","")"
$ws.Range("H23").Value = "skos:exactMatch(separator="","")"
$ws.Range("I23").Value = "skos:closeMatch(separator="","")"
$ws.Range("J23").Value = "skos:broadMatch(separator="","")"
$ws.Range("K23").Value = "iop:hasProperty"
$ws.Range("L23").Value = "iop:hasObjectOfInterest"
$ws.Range("M23").Value = "iop:hasMatrix"
$ws.Range("N23").Value = "iop:hasContextObject(separator="","")"
$ws.Range("O23").Value = "iop:hasConstraint(separator="","")"
$ws.Range("P23").Value = "puv:statistic(separator="","")"
$ws.Range("Q23").Value = "puv:usesMethod(separator="","")"
$ws.Range("R23").Value = "sosa:madeBySensor(separator="","")"
$ws.Range("S23").Value = "puv:uom(separator="","")"
$ws.Range("T23").Value = "owl:deprecated^^xsd:boolean"
$ws.Range("U23").Value = "skos:editorialNote@en"
$ws.Range("V23").Value = "rdf:type"
$ws.Range("W23").Value = "dct:modified^^xsd:date"
$ws.Range("X23").Value = "dct:created^^xsd:date"
$ws.Range("Y23").Value = "dct:creator(separator="","")"
$ws.Range("Z23").Value = "dct:contributor(separator="","")"
$ws.Range("C24").ClearContents()
$ws.Range("C27").Value = "new"
$ws.Range("E27").ClearContents()
$ws.Range("G27").Value = "vocab:1002"
$ws.Range("C28").Value = "intermediate"
$ws.Range("E28").ClearContents()
$ws.Range("G28").Value = "vocab:1002"
$ws.Range("E29").ClearContents()
$ws.Range("G29").Value = "vocab:1002"
$ws.Range("F30").ClearContents()
$ws.Range("V30").Value = "owl:ObjectProperty"

$ws.Range("AP1").EntireColumn.Delete()
